# Actualización 10 de Mayo
# Updates the statistics sheets (1P, 2P, Final) with new figures and
# clears out the "Rescatables" (make-up exam) list, which is now empty.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Estadisticos 1P"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Estadisticos 1P")

$ws1.Range("D2").Value2 = 1
$ws1.Range("E2").Value2 = 5
$ws1.Range("F2").Value2 = 22
$ws1.Range("G2").Value2 = 78.57
$ws1.Range("H2").Value2 = 6.6

$ws1.Range("D3").Value2 = 4
$ws1.Range("E3").Value2 = 3
$ws1.Range("F3").Value2 = 21
$ws1.Range("G3").Value2 = 75
$ws1.Range("H3").Value2 = 6.3

$ws1.Range("D4").Value2 = 3
$ws1.Range("E4").Value2 = 3
$ws1.Range("F4").Value2 = 16
$ws1.Range("G4").Value2 = 72.73
$ws1.Range("H4").Value2 = 6.4

# ---------------------------------------------------------------------
# Sheet "Estadisticos 2P"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Estadisticos 2P")

$ws2.Range("D2").Value2 = 11
$ws2.Range("E2").Value2 = 10
$ws2.Range("F2").Value2 = 17
$ws2.Range("G2").Value2 = 60.71
$ws2.Range("H2").Value2 = 7.2

$ws2.Range("D3").Value2 = 13
$ws2.Range("E3").Value2 = 9
$ws2.Range("F3").Value2 = 15
$ws2.Range("G3").Value2 = 53.57
$ws2.Range("H3").Value2 = 6.5

$ws2.Range("D4").Value2 = 4
$ws2.Range("E4").Value2 = 1
$ws2.Range("F4").Value2 = 18
$ws2.Range("G4").Value2 = 81.82
$ws2.Range("H4").Value2 = 6.4

# ---------------------------------------------------------------------
# Sheet "Estadisticos Final"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Estadisticos Final")

$ws3.Range("D2").Value2 = 1
$ws3.Range("E2").Value2 = 5
$ws3.Range("F2").Value2 = 22
$ws3.Range("G2").Value2 = 78.57
$ws3.Range("H2").Value2 = 6.8

$ws3.Range("D3").Value2 = 4
$ws3.Range("E3").Value2 = 3
$ws3.Range("F3").Value2 = 21
$ws3.Range("G3").Value2 = 75

$ws3.Range("D4").Value2 = 3
$ws3.Range("F4").Value2 = 19
$ws3.Range("G4").Value2 = 86.36

# ---------------------------------------------------------------------
# Sheet "Rescatables" - remove the student rows, only header remains
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Rescatables")
$ws4.Rows("2:6").Delete()
